# Add a new Translation row (row 13) for a new UI text entry.
# Mirrors the pattern of the existing rows in the Translation sheet:
#   B = Text ID, C = Typography Name, D = Alignment, E = Translated text (GB), F = Direction

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Translation")

$ws.Cells.Item(13, 2).Value = "SingleUseId13"
$ws.Cells.Item(13, 3).Value = "OurTypography"
$ws.Cells.Item(13, 4).Value = "Center"
$ws.Cells.Item(13, 5).Value = "FASTER WITH THAT CARD`nYOU HOE FUCKING"
$ws.Cells.Item(13, 6).Value = "LTR"

# Keep the new cells using the default/unstyled formatting (matching the
# rest of the data rows) and make sure the row height stays at the sheet's
# standard height instead of the multi-line auto height.
$ws.Range("B13:F13").Style = "Normal"
$ws.Rows.Item(13).AutoFit()
